{"js": "// \"Version 1.\" -> \"Version 2.\" (wireframes.docx version bump)\nconst body = context.document.body;\n\n// Only touch the version digit itself (\" 1.\" -> \" 2.\") so the rest of the\n// paragraph (the word \"Version\", the bookmark, etc.) is left untouched -\n// this is the minimal edit that matches what actually changed.\nconst results = body.search(\" 1.\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\" 2.\", Word.InsertLocation.replace);\n} else {\n  // Fallback in case the leading space isn't part of the match for some\n  // reason - fall back to replacing the whole \"Version 1.\" phrase.\n  const whole = body.search(\"Version 1.\", { matchCase: true });\n  whole.load(\"text\");\n  await context.sync();\n  if (whole.items.length > 0) {\n    whole.items[0].insertText(\"Version 2.\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Version 1.\" -> \"Version 2.\" (wireframes.docx version bump)\n$d = $word.ActiveDocument\n\n# Only touch the version digit itself (\" 1.\" -> \" 2.\") so the rest of the\n# paragraph (the word \"Version\", the bookmark, etc.) is left untouched -\n# this is the minimal edit that matches what actually changed.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" 1.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \" 2.\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback in case the leading space isn't part of the match for some\n    # reason - fall back to replacing the whole \"Version 1.\" phrase.\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Text = \"Version 1.\"\n    $find2.Replacement.ClearFormatting()\n    $find2.Replacement.Text = \"Version 2.\"\n    $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n}\n"}
